$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A8").Value = "03/05/2021 Taller Parte movil"
